# Auto-generated edit script applying the Midgardsormr_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 884.56525
$ws.Range("J17").Value = 884.56525
$ws.Range("L17").Value = 2653.69575
$ws.Range("N17").Value = -2989.69575

# Row 86
$ws.Range("H86").Value = 1866.8695
$ws.Range("I86").Value = 1623.1428
$ws.Range("K86").Value = 1623.1428
$ws.Range("M86").Value = -500.1428000000001

# Row 89
$ws.Range("H89").Value = 1866.8695
$ws.Range("I89").Value = 1623.1428
$ws.Range("K89").Value = 8115.714
$ws.Range("M89").Value = -2499.714

# Row 106
$ws.Range("H106").Value = 3232.111
$ws.Range("I106").Value = 3232.111
$ws.Range("K106").Value = 3232.111
$ws.Range("M106").Value = -2601.111

# Row 113
$ws.Range("H113").Value = 5066.5557
$ws.Range("J113").Value = 6399
$ws.Range("L113").Value = 6399
$ws.Range("N113").Value = -12907

# Row 128
$ws.Range("H128").Value = 115000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 115000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 115000
$ws.Range("N128").Value = -124960
$ws.Range("M128").ClearContents()

# Row 132
$ws.Range("H132").Value = 17233.422
$ws.Range("I132").Value = 18598.39
$ws.Range("K132").Value = 55795.17
$ws.Range("M132").Value = -53265.17

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4382.9287
$ws.Range("J45").Value = 5712.4287
$ws.Range("L45").Value = 5712.4287
$ws.Range("N45").Value = -6466.4287

# Row 98
$ws.Range("H98").Value = 39332.668
$ws.Range("J98").Value = 39332.668
$ws.Range("L98").Value = 39332.668
$ws.Range("N98").Value = -45322.668

# Row 132
$ws.Range("H132").Value = 1629.6842
$ws.Range("I132").Value = 1411.1613
$ws.Range("J132").Value = 2597.4285
$ws.Range("K132").Value = 4233.4839
$ws.Range("L132").Value = 7792.2855
$ws.Range("M132").Value = -1703.4839
$ws.Range("N132").Value = -12852.2855

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 49996.5
$ws.Range("J81").Value = 49996.5
$ws.Range("L81").Value = 49996.5
$ws.Range("N81").Value = -52118.5

# Row 84
$ws.Range("H84").Value = 49996.5
$ws.Range("J84").Value = 49996.5
$ws.Range("L84").Value = 149989.5
$ws.Range("N84").Value = -160597.5

# Row 99
$ws.Range("H99").Value = 2082.4707
$ws.Range("I99").Value = 2042.6666
$ws.Range("K99").Value = 2042.6666
$ws.Range("M99").Value = -544.6666

# Row 105
$ws.Range("H105").Value = 3196.5625
$ws.Range("I105").Value = 3018.1538
$ws.Range("K105").Value = 3018.1538
$ws.Range("M105").Value = -1271.1538

# Row 107
$ws.Range("H107").Value = 5167.96
$ws.Range("J107").Value = 4382.7144
$ws.Range("L107").Value = 4382.7144
$ws.Range("N107").Value = -8222.714400000001

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 56883.723
$ws.Range("I132").Value = 67513.87
$ws.Range("J132").Value = 3733
$ws.Range("K132").Value = 202541.61
$ws.Range("L132").Value = 11199
$ws.Range("M132").Value = -200011.61
$ws.Range("N132").Value = -16259

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 4477.086
$ws.Range("I68").Value = 1066
$ws.Range("K68").Value = 3198
$ws.Range("M68").Value = -2387

# Row 71
$ws.Range("H71").Value = 4477.086
$ws.Range("I71").Value = 1066
$ws.Range("K71").Value = 9594
$ws.Range("M71").Value = -5538

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 19
$ws.Range("H19").Value = 5004710
$ws.Range("J19").Value = 5004710
$ws.Range("L19").Value = 5004710
$ws.Range("N19").Value = -5005286

# Row 93
$ws.Range("H93").Value = 79999
$ws.Range("J93").Value = 79999
$ws.Range("L93").Value = 79999
$ws.Range("N93").Value = -83743

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1823.3125
$ws.Range("I22").Value = 1417.8889
$ws.Range("J22").Value = 2344.5715
$ws.Range("K22").Value = 1417.8889
$ws.Range("L22").Value = 2344.5715
$ws.Range("M22").Value = -1122.8889
$ws.Range("N22").Value = -2934.5715

# Row 27
$ws.Range("H27").Value = 1823.3125
$ws.Range("I27").Value = 1417.8889
$ws.Range("J27").Value = 2344.5715
$ws.Range("K27").Value = 1417.8889
$ws.Range("L27").Value = 2344.5715
$ws.Range("M27").Value = -1310.8889
$ws.Range("N27").Value = -2558.5715

# Row 110
$ws.Range("H110").Value = 79999.5
$ws.Range("J110").Value = 79999.5
$ws.Range("L110").Value = 79999.5
$ws.Range("N110").Value = -88179.5

# Row 136
$ws.Range("H136").Value = 3458.4736
$ws.Range("I136").Value = 3209.258
$ws.Range("K136").Value = 9627.773999999999
$ws.Range("M136").Value = -7077.773999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# Row 62
$ws.Range("H62").Value = 2639.5833
$ws.Range("I62").Value = 2922.4443
$ws.Range("J62").Value = 1791
$ws.Range("K62").Value = 2922.4443
$ws.Range("L62").Value = 1791
$ws.Range("M62").Value = -2298.4443
$ws.Range("N62").Value = -3039

# Row 65
$ws.Range("H65").Value = 2639.5833
$ws.Range("I65").Value = 2922.4443
$ws.Range("J65").Value = 1791
$ws.Range("K65").Value = 14612.2215
$ws.Range("L65").Value = 8955
$ws.Range("M65").Value = -11492.2215
$ws.Range("N65").Value = -15195

# Row 81
$ws.Range("H81").Value = 7798.5
$ws.Range("I81").Value = 11124.875
$ws.Range("J81").Value = 4472.125
$ws.Range("K81").Value = 22249.75
$ws.Range("L81").Value = 8944.25
$ws.Range("M81").Value = -21188.75
$ws.Range("N81").Value = -11066.25

# Row 84
$ws.Range("H84").Value = 7798.5
$ws.Range("I84").Value = 11124.875
$ws.Range("J84").Value = 4472.125
$ws.Range("K84").Value = 111248.75
$ws.Range("L84").Value = 44721.25
$ws.Range("M84").Value = -105944.75
$ws.Range("N84").Value = -55329.25

# Row 132
$ws.Range("H132").Value = 1314.079
$ws.Range("I132").Value = 950.2759
$ws.Range("K132").Value = 2850.8277
$ws.Range("M132").Value = -320.8276999999998
